# Backlog.xlsx edit: "début documentation livrable 2"
#
# Renumbers/re-prioritises the product backlog: several user stories
# (admin rights levels, super-admin registration/rights management,
# item reservation) are dropped, the "admin can edit sections" story is
# simplified (no more super-admin distinction) and moved to the bottom
# of the currently-scoped set, and both the Backlog sheet and the
# Sprint 1 sheet are trimmed down to match.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Backlog")
$ws2 = $wb.Worksheets.Item("Sprint 1")

# ---------------------------------------------------------------------
# Final row layout shared by both sheets for rows 2-8 (identical text)
# ---------------------------------------------------------------------
function Set-CommonRows($ws) {
    $ws.Range("A2").Value = 1
    $ws.Range("C2").Value = 75
    $ws.Range("D2").Value = 1
    $ws.Range("E2").Value = "En tant que visiteur, je souhaite qu'il y ai une page d'accueil afin de me permettre de naviguer dans le site."

    $ws.Range("B3").Value = "a"
    $ws.Range("E3").Value = "Le site web a un template pour la page d'accueil avec le menu"

    $ws.Range("A4").Value = 2
    $ws.Range("C4").Value = 4
    $ws.Range("D4").Value = 2
    $ws.Range("E4").Value = "En tant que administrateur, je souhaite pouvoir me connecter  afin de pouvoir bénéficier de mes droits d'administrateurs."

    $ws.Range("B5").Value = "a"
    $ws.Range("E5").Value = "La connection se fait avec un email et un mot de passe."

    $ws.Range("B6").Value = "b"
    $ws.Range("E6").Value = "Il devient impossible d'effectuer une connexion si l'administrateur est déjà connecté."

    $ws.Range("A7").Value = 3
    $ws.Range("C7").Value = 4
    $ws.Range("D7").Value = 3
    $ws.Range("E7").Value = "En tant que administrateur, je souhaite pouvoir me déconnecter afin de m'assurer que personne n'utilise ma session."

    $ws.Range("B8").Value = "a"
    $ws.Range("E8").Value = "La déconnexion ramène à la page d'accueil."
}

Set-CommonRows $ws1
Set-CommonRows $ws2

# ---------------------------------------------------------------------
# Backlog sheet only: rows 9-21 (stories kept further down the backlog)
# ---------------------------------------------------------------------
$ws1.Range("A9").Value = 4
$ws1.Range("C9").Value = 50
$ws1.Range("D9").Value = 4
$ws1.Range("E9").Value = "En tant que administrateur, je souhaire pouvoir modifier les sections du sites afin de maintenir le site à jour."

$ws1.Range("B10").Value = "a"
$ws1.Range("E10").Value = "Seul les administrateur peuvent modifier les section du site"

$ws1.Range("A11").Value = 5
$ws1.Range("C11").Value = 6
$ws1.Range("D11").Value = 5
$ws1.Range("E11").Value = "En tant que visiteur, je souhaite pouvoir effectuer un don afin de financer les paroisses."

$ws1.Range("B12").Value = "a"
$ws1.Range("E12").Value = "Il y a un bouton de donation visible sur toute les pages."

# Row 13 used to hold a full story (A/C/D/E) - it becomes a plain
# criteria row now, so clear the leftover A/C/D cells first.
$ws1.Range("A13:D13").ClearContents()
$ws1.Range("B13").Value = "b"
$ws1.Range("E13").Value = "Le bouton de donation redirige vers PayPal."

# Row 14 used to be a plain criteria row - it becomes a full story row.
$ws1.Range("B14").ClearContents()
$ws1.Range("A14").Value = 6
$ws1.Range("C14").Value = 30
$ws1.Range("D14").Value = 6
$ws1.Range("E14").Value = "En tant que visiteur, je souhaite pouvoir me renseigner sur les paroisses afin de me tenir au courant des nouveauté de celle-ci"

# Row 15 used to be a full story row - it becomes a plain criteria row.
$ws1.Range("A15:D15").ClearContents()
$ws1.Range("B15").Value = "a"
$ws1.Range("E15").Value = "Les nouvelles sont facilement visible sur le site."

$ws1.Range("B16").Value = "b"
$ws1.Range("E16").Value = "Les nouvelles sont groupés par paroisse."

# Row 17 used to be a full story row (A/C/D) - now plain criteria, clear it.
$ws1.Range("A17:D17").ClearContents()
$ws1.Range("B17").Value = "c"
$ws1.Range("E17").Value = "Les nouvelles s'affiches durant une période spécifique."

# Row 18 used to be a plain criteria row - it becomes a full story row.
$ws1.Range("B18").ClearContents()
$ws1.Range("A18").Value = 7
$ws1.Range("C18").Value = 25
$ws1.Range("D18").Value = 7
$ws1.Range("E18").Value = "En tant que visiteur, je souhaite pouvoir remplir des formulaires afin de m'inscrire aux différentes communautés religieuse"

$ws1.Range("B19").Value = "a"
$ws1.Range("E19").Value = "Les informations du visiteurs dans les formulaires sont encryptés. "

$ws1.Range("A20").Value = 8
$ws1.Range("C20").Value = 6
$ws1.Range("D20").Value = 8
$ws1.Range("E20").Value = "En tant que visiteur, je souhaite pouvoir m'inscire comme bénévole afin de  contribuer à la communauté."

$ws1.Range("B21").Value = "a"
$ws1.Range("E21").Value = "L'inscription comme bénévole se fait par formulaire."

# Drop everything past row 21 on the Backlog sheet.
$ws1.Rows("22:35").Delete()

# Drop the trailing two rows on the Sprint 1 sheet (used to hold the
# "register a new admin" story, which is gone from scope entirely).
$ws2.Rows("9:10").Delete()

# ---------------------------------------------------------------------
# Selections / active sheet: Backlog becomes the active (first) tab.
# ---------------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("A7").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("F9").Select() | Out-Null
